$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel must be pre-formatted as Text
# so they are stored as strings (matching the source data export format),
# not auto-converted to numbers.
$textCells = @("D5", "D6", "D9", "D11", "D12", "D13", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D42", "D44", "D45", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '71.049.38'
$ws.Range("E2").Value = '  +6.53%  '
$ws.Range("D3").Value = '3.683.76'
$ws.Range("E3").Value = '  +18.99%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '598.97'
$ws.Range("E5").Value = '  +3.96%  '
$ws.Range("D6").Value = '184.62'
$ws.Range("E6").Value = '  +6.76%  '
$ws.Range("D7").Value = '3.681.26'
$ws.Range("E7").Value = '  +18.92%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.537'
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("E10").Value = '  +7.82%  '
$ws.Range("D11").Value = '6.62'
$ws.Range("E11").Value = '  +3.91%  '
$ws.Range("D12").Value = '0.500'
$ws.Range("E12").Value = '  +5.51%  '
$ws.Range("D13").Value = '40.18'
$ws.Range("E13").Value = '  +12.34%  '
$ws.Range("E14").Value = '  +6.15%  '
$ws.Range("D15").Value = '4.296.81'
$ws.Range("E15").Value = '  +19.03%  '
$ws.Range("D16").Value = '71.119.39'
$ws.Range("E16").Value = '  +6.65%  '
$ws.Range("D17").Value = '3.676.42'
$ws.Range("E17").Value = '  +18.84%  '
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("E19").Value = '  +7.58%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '514.70'
$ws.Range("E21").Value = '  +6.58%  '
$ws.Range("D22").Value = '9.21'
$ws.Range("E22").Value = '  +18.23%  '
$ws.Range("E23").Value = '  +7.83%  '
$ws.Range("D24").Value = '87.59'
$ws.Range("E24").Value = '  +5.10%  '
$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +11.56%  '
$ws.Range("D26").Value = '13.53'
$ws.Range("E26").Value = '  +7.10%  '
$ws.Range("E27").Value = '  +10.52%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '2.53'
$ws.Range("E29").Value = '  +11.88%  '
$ws.Range("D30").Value = '8.19'
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0000111'
$ws.Range("E31").Value = '  +19.26%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '2.78'
$ws.Range("E32").Value = '  +7.29%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '31.65'
$ws.Range("E33").Value = '  +13.20%  '
$ws.Range("D34").Value = '0.117'
$ws.Range("E34").Value = '  +4.63%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '6.12'
$ws.Range("E36").Value = '  +9.49%  '
$ws.Range("E37").Value = '  +7.89%  '
$ws.Range("E38").Value = '  +12.14%  '
$ws.Range("D39").Value = '2.17'
$ws.Range("E39").Value = '  +10.46%  '
$ws.Range("D40").Value = '51.10'
$ws.Range("E40").Value = '  +4.26%  '
$ws.Range("E41").Value = '  +4.31%  '
$ws.Range("D42").Value = '45.52'
$ws.Range("E42").Value = '  -5.13%  '
$ws.Range("D43").Value = '3.147.83'
$ws.Range("E43").Value = '  +12.50%  '
$ws.Range("D44").Value = '8.85'
$ws.Range("E44").Value = '  +7.03%  '
$ws.Range("D45").Value = '413.96'
$ws.Range("E45").Value = '  +12.24%  '
$ws.Range("E46").Value = '  +5.71%  '
$ws.Range("E47").Value = '  +6.54%  '
$ws.Range("D48").Value = '28.30'
$ws.Range("E48").Value = '  +15.91%  '
$ws.Range("D49").Value = '137.74'
$ws.Range("E49").Value = '  +2.57%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").Value = '2.47'
$ws.Range("E51").Value = '  +12.75%  '
